$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.782.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.98%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.333.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.43"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.82"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.44%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.53%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.92"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.08"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0799"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.54%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.80"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.79%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.392.00"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.14%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.724.17"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.85%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.24"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.90%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.61"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -7.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.69"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.26"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.99"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.56"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.21%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.99"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.33"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.77"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -8.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.36"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.50"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.73%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.11"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.34"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0726"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.95"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.62%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.01%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.02%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.35"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.022.96"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0284"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.60"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -7.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.28"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.92"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.02"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.92"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.559.49"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.46%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.77%  "

# Row 34 and Row 36 swap (RenderToken <-> WEMIXToken)
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.45"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.35%  "

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.58"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.39%  "
